$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1089.1428
$ws.Range("I2").Value = 1062.4166
$ws.Range("K2").Value = 1062.4166
$ws.Range("M2").Value = -949.4166
$ws.Range("H17").Value = 5365.0615
$ws.Range("J17").Value = 5684.754
$ws.Range("L17").Value = 17054.262
$ws.Range("N17").Value = -17390.262
$ws.Range("H21").Value = 1999.5
$ws.Range("I21").Value = 1999.5
$ws.Range("K21").Value = 1999.5
$ws.Range("M21").Value = -1531.5
$ws.Range("H23").Value = 1999.5
$ws.Range("I23").Value = 1999.5
$ws.Range("K23").Value = 1999.5
$ws.Range("M23").Value = -1765.5
$ws.Range("H53").Value = 604.9583
$ws.Range("I53").Value = 232.14285
$ws.Range("J53").Value = 1126.9
$ws.Range("K53").Value = 232.14285
$ws.Range("L53").Value = 1126.9
$ws.Range("M53").Value = 404.85715
$ws.Range("N53").Value = -2400.9
$ws.Range("H87").Value = 89500
$ws.Range("J87").Value = 89500
$ws.Range("L87").Value = 89500
$ws.Range("N87").Value = -91996
$ws.Range("H90").Value = 89500
$ws.Range("J90").Value = 89500
$ws.Range("L90").Value = 268500
$ws.Range("N90").Value = -280980
$ws.Range("H128").Value = 100000
$ws.Range("J128").Value = 100000
$ws.Range("L128").Value = 100000
$ws.Range("N128").Value = -109960
$ws.Range("H137").Value = 180295.5
$ws.Range("I137").Value = 445569.5
$ws.Range("K137").Value = 1336708.5
$ws.Range("M137").Value = -1334158.5
$ws.Range("H138").Value = 2575.7593
$ws.Range("I138").Value = 1135.8788
$ws.Range("J138").Value = 4838.4287
$ws.Range("K138").Value = 3407.6364
$ws.Range("L138").Value = 14515.2861
$ws.Range("M138").Value = 1732.3636
$ws.Range("N138").Value = -24795.2861

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4971.7754
$ws.Range("I32").Value = 2531
$ws.Range("K32").Value = 2531
$ws.Range("M32").Value = -2244
$ws.Range("H45").Value = 10994069
$ws.Range("I45").Value = 17097222
$ws.Range("K45").Value = 17097222
$ws.Range("M45").Value = -17096845
$ws.Range("H61").Value = 2780.9333
$ws.Range("I61").Value = 2243.25
$ws.Range("J61").Value = 4931.6665
$ws.Range("K61").Value = 2243.25
$ws.Range("L61").Value = 4931.6665
$ws.Range("M61").Value = -2031.25
$ws.Range("N61").Value = -5355.6665
$ws.Range("H74").Value = 95235.91
$ws.Range("I74").Value = 10121.111
$ws.Range("K74").Value = 10121.111
$ws.Range("M74").Value = -9247.111000000001
$ws.Range("H77").Value = 95235.91
$ws.Range("I77").Value = 10121.111
$ws.Range("K77").Value = 50605.55500000001
$ws.Range("M77").Value = -46237.55500000001
$ws.Range("H110").Value = 1209527.9
$ws.Range("I110").Value = 1635712.6
$ws.Range("K110").Value = 1635712.6
$ws.Range("M110").Value = -1633667.6
$ws.Range("H132").Value = 2264.111
$ws.Range("I132").Value = 1941.3334
$ws.Range("K132").Value = 5824.0002
$ws.Range("M132").Value = -3294.0002
$ws.Range("H136").Value = 2780.9333
$ws.Range("I136").Value = 2243.25
$ws.Range("J136").Value = 4931.6665
$ws.Range("K136").Value = 6729.75
$ws.Range("L136").Value = 14794.9995
$ws.Range("M136").Value = -4179.75
$ws.Range("N136").Value = -19894.9995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2749178.5
$ws.Range("I107").Value = 3761093.8
$ws.Range("J107").Value = 2551.5715
$ws.Range("K107").Value = 3761093.8
$ws.Range("L107").Value = 2551.5715
$ws.Range("M107").Value = -3759173.8
$ws.Range("N107").Value = -6391.5715
$ws.Range("H134").Value = 3265.1428
$ws.Range("I134").Value = 1590.0322
$ws.Range("K134").Value = 4770.096600000001
$ws.Range("M134").Value = -2235.096600000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2081.5
$ws.Range("I16").Value = 1773.5
$ws.Range("K16").Value = 1773.5
$ws.Range("M16").Value = -1486.5
$ws.Range("H19").Value = 2857506
$ws.Range("I19").Value = 3333586.2
$ws.Range("J19").Value = 1025
$ws.Range("K19").Value = 3333586.2
$ws.Range("L19").Value = 1025
$ws.Range("M19").Value = -3333416.2
$ws.Range("N19").Value = -1365
$ws.Range("H24").Value = 2857506
$ws.Range("I24").Value = 3333586.2
$ws.Range("J24").Value = 1025
$ws.Range("K24").Value = 3333586.2
$ws.Range("L24").Value = 1025
$ws.Range("M24").Value = -3333416.2
$ws.Range("N24").Value = -1365
$ws.Range("H31").Value = 4179.23
$ws.Range("I31").Value = 1337.5
$ws.Range("K31").Value = 1337.5
$ws.Range("M31").Value = -1042.5
$ws.Range("H34").Value = 4179.23
$ws.Range("I34").Value = 1337.5
$ws.Range("K34").Value = 1337.5
$ws.Range("M34").Value = -1135.5
$ws.Range("H50").Value = 5416.5835
$ws.Range("J50").Value = 5416.5835
$ws.Range("L50").Value = 5416.5835
$ws.Range("N50").Value = -6666.5835
$ws.Range("H58").Value = 2577.2856
$ws.Range("I58").Value = 1756.8334
$ws.Range("J58").Value = 3671.2222
$ws.Range("K58").Value = 1756.8334
$ws.Range("L58").Value = 3671.2222
$ws.Range("M58").Value = -1553.8334
$ws.Range("N58").Value = -4077.2222
$ws.Range("H68").Value = 57272.363
$ws.Range("J68").Value = 65555.22
$ws.Range("L68").Value = 65555.22
$ws.Range("N68").Value = -67053.22
$ws.Range("H71").Value = 57272.363
$ws.Range("J71").Value = 65555.22
$ws.Range("L71").Value = 196665.66
$ws.Range("N71").Value = -204153.66
$ws.Range("H113").Value = 2081.5
$ws.Range("I113").Value = 1773.5
$ws.Range("K113").Value = 1773.5
$ws.Range("M113").Value = 396.5
$ws.Range("H122").Value = 3640.3
$ws.Range("I122").Value = 3337.3333
$ws.Range("K122").Value = 10011.9999
$ws.Range("M122").Value = -7561.999899999999
$ws.Range("H134").Value = 3443.2415
$ws.Range("I134").Value = 3242.0625
$ws.Range("J134").Value = 3690.8462
$ws.Range("K134").Value = 9726.1875
$ws.Range("L134").Value = 11072.5386
$ws.Range("M134").Value = -7191.1875
$ws.Range("N134").Value = -16142.5386
$ws.Range("H135").Value = 118799.336
$ws.Range("J135").Value = 118799.336
$ws.Range("L135").Value = 118799.336
$ws.Range("N135").Value = -128939.336
$ws.Range("H136").Value = 2577.2856
$ws.Range("I136").Value = 1756.8334
$ws.Range("J136").Value = 3671.2222
$ws.Range("K136").Value = 5270.5002
$ws.Range("L136").Value = 11013.6666
$ws.Range("M136").Value = -2720.5002
$ws.Range("N136").Value = -16113.6666

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4134856
$ws.Range("J4").Value = 267980
$ws.Range("L4").Value = 803940
$ws.Range("N4").Value = -804164
$ws.Range("H44").Value = 500200
$ws.Range("I44").Value = 400
$ws.Range("J44").Value = 1000000
$ws.Range("K44").Value = 1200
$ws.Range("L44").Value = 3000000
$ws.Range("M44").Value = -802
$ws.Range("N44").Value = -3000796
$ws.Range("H98").Value = 596.4545000000001
$ws.Range("I98").Value = 545.6667
$ws.Range("J98").Value = 615.5
$ws.Range("K98").Value = 1637.0001
$ws.Range("L98").Value = 1846.5
$ws.Range("M98").Value = -139.0001
$ws.Range("N98").Value = -4842.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 27779346
$ws.Range("I113").Value = 33334814
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 33334814
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -33332644
$ws.Range("N113").Value = -6340

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H16").Value = 510.72726
$ws.Range("I16").Value = 526
$ws.Range("J16").Value = 190
$ws.Range("K16").Value = 526
$ws.Range("L16").Value = 190
$ws.Range("M16").Value = -356
$ws.Range("N16").Value = -530
$ws.Range("H61").Value = 12351712
$ws.Range("I61").Value = 15879272
$ws.Range("J61").Value = 5250
$ws.Range("K61").Value = 15879272
$ws.Range("L61").Value = 5250
$ws.Range("M61").Value = -15879070
$ws.Range("N61").Value = -5654
$ws.Range("H113").Value = 12351712
$ws.Range("I113").Value = 15879272
$ws.Range("J113").Value = 5250
$ws.Range("K113").Value = 15879272
$ws.Range("L113").Value = 5250
$ws.Range("M113").Value = -15877102
$ws.Range("N113").Value = -9590
$ws.Range("H136").Value = 30796.305
$ws.Range("I136").Value = 43719.457
$ws.Range("K136").Value = 131158.371
$ws.Range("M136").Value = -128608.371

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 111119060
$ws.Range("I107").Value = 125005180
$ws.Range("K107").Value = 375015540
$ws.Range("M107").Value = -375013620
$ws.Range("H109").Value = 69989.5
$ws.Range("J109").Value = 69989.5
$ws.Range("L109").Value = 69989.5
$ws.Range("N109").Value = -72763.5
$ws.Range("H113").Value = 1024.2727
$ws.Range("I113").Value = 896.6667
$ws.Range("J113").Value = 1072.125
$ws.Range("K113").Value = 2690.0001
$ws.Range("L113").Value = 3216.375
$ws.Range("M113").Value = -520.0001000000002
$ws.Range("N113").Value = -7556.375
$ws.Range("H126").Value = 2769.6924
$ws.Range("I126").Value = 2640.48
$ws.Range("K126").Value = 7921.440000000001
$ws.Range("M126").Value = -5451.440000000001
$ws.Range("H136").Value = 2460.2
$ws.Range("J136").Value = 5899.8
$ws.Range("L136").Value = 17699.4
$ws.Range("N136").Value = -22799.4
